# Logic tree input file updated
#
# The decision tree gets two new "default answer" rows (Node2 =
# "Possible_Problem", Relationship = the generic possible-problem blurb),
# one inserted right under each of the two "...freeway/highway driving or
# ... traffic?" question rows - mirroring the pattern already used for the
# two earlier question rows (current rows 4 and 7). Every row below each
# insertion point shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$possibleProblem = "Possible_Problem"
$defaultAnswer = "Possible_Problem:25% Front wheel speed sensor`n20% Front Hub Assembly`n15% Rear Wheel Speed Sensor`n15% Speed Sensor Tone Ring`n10% Brake Pedal Sensor / Stop Light Switch`n10% Rusted Speed Sensor Mounting Surface`n5% ABS Control Unit"

$noiseQuestion = "Problem:Is the noise present during freeway/highway driving or Stop & Go? (Please answer as: Stop & Go, Freeway, Both)"
$lightQuestion = "Problem:Does the light come on during freeway/highway driving or City traffic? (Please answer as: Stop & Go, Freeway, Both)"

# --- Insert new row 10 (under the "noise" freeway/Stop & Go question, which
#     starts at the existing row 10) ---
$ws.Rows.Item(10).Insert()

$ws.Cells.Item(10, 1).Value = $noiseQuestion
$ws.Cells.Item(10, 2).Value = $possibleProblem
$ws.Cells.Item(10, 3).Value = $defaultAnswer
$ws.Cells.Item(10, 3).WrapText = $true
$ws.Rows.Item(10).RowHeight = 409.6

# --- Insert new row 14 (under the "light" freeway/City traffic question).
#     After the row-10 insertion above, the old row 13 (the first row of the
#     "light" question block) is now row 14. ---
$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14, 1).Value = $lightQuestion
$ws.Cells.Item(14, 2).Value = $possibleProblem
$ws.Cells.Item(14, 3).Value = $defaultAnswer
$ws.Cells.Item(14, 3).WrapText = $true
$ws.Rows.Item(14).RowHeight = 409.6

# --- Update the saved view state to match: scrolled down with A18 selected ---
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("A18").Select()
